# Natmi following Dr Hou advice
# Update the LR-pair result values in rows 2-7 (columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        E = 3; G = 1.388571333333333; H = 4.165713999999999
        I = 0.3523526610542377; J = 0.3523526610542378
        K = 3
        M = 125.0114046666667; N = 375.034214
        O = 0.4125441987306753; P = 0.4125441987306753
        Q = 173.5872528598662; R = 1562.285275738796
        S = 0.1453610462252417; T = 0.1453610462252417
    }
    3 = @{
        E = 3; G = 1.388571333333333; H = 4.165713999999999
        I = 0.3523526610542377; J = 0.3523526610542378
        K = 3
        M = 139.2310486666667; N = 417.693146
        O = 0.4594697704883666; P = 0.4594697704883666
        Q = 193.3322428884715; R = 1739.990185996244
        S = 0.1618953963055558; T = 0.1618953963055559
    }
    4 = @{
        E = 3; G = 1.388571333333333; H = 4.165713999999999
        I = 0.3523526610542377; J = 0.3523526610542378
        K = 3
        M = 38.78302866666667; N = 116.349086
        O = 0.1279860307809581; P = 0.1279860307809581
        Q = 53.85300182637822; R = 484.6770164374039
        S = 0.04509621852344017; T = 0.04509621852344017
    }
    5 = @{
        E = 3; G = 2.552285333333333; H = 7.656856
        I = 0.6476473389457622; J = 0.6476473389457623
        K = 3
        M = 125.0114046666667; N = 375.034214
        O = 0.4125441987306753; P = 0.4125441987306753
        Q = 319.0647746301316; R = 2871.582971671184
        S = 0.2671831525054336; T = 0.2671831525054336
    }
    6 = @{
        E = 3; G = 2.552285333333333; H = 7.656856
        I = 0.6476473389457622; J = 0.6476473389457623
        K = 3
        M = 139.2310486666667; N = 417.693146
        O = 0.4594697704883666; P = 0.4594697704883666
        Q = 355.3573634565528; R = 3198.216271108976
        S = 0.2975743741828107; T = 0.2975743741828107
    }
    7 = @{
        E = 3; G = 2.552285333333333; H = 7.656856
        I = 0.6476473389457622; J = 0.6476473389457623
        K = 3
        M = 38.78302866666667; N = 116.349086
        O = 0.1279860307809581; P = 0.1279860307809581
        Q = 98.98535524817956; R = 890.868197233616
        S = 0.08288981225751793; T = 0.08288981225751794
    }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
